$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$variavel = "Número médio de moradores"

$regioes = @("Brasil", "Nordeste", "Sergipe")
$anos = @("31/12/2016", "31/12/2017", "31/12/2018", "31/12/2019", "31/12/2022", "31/12/2023")
$valores = @{
    "Brasil"   = @(3, 3, 3, 3, 2.9, 2.8)
    "Nordeste" = @(3.2, 3.2, 3.1, 3.1, 3, 2.8)
    "Sergipe"  = @(3, 3, 3.1, 3, 2.9, 2.8)
}

$row = 2
foreach ($regiao in $regioes) {
    $vals = $valores[$regiao]
    for ($i = 0; $i -lt $anos.Length; $i++) {
        $ws.Cells.Item($row, 1).Value = $regiao
        $ws.Cells.Item($row, 2).Value = $variavel
        $ws.Cells.Item($row, 3).Value = $anos[$i]
        $ws.Cells.Item($row, 4).Value = $vals[$i]
        $row = $row + 1
    }
}
